$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.927.00'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.89%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.632.64'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.503'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.47%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.69'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0788'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.856.52'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.599.33'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.23%  '
$ws.Range("E15").Value = '  -1.73%  '
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.892.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.70%  '
$ws.Range("E20").Value = '  -1.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '193.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.96'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.26'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.47%  '
$ws.Range("E24").Value = '  -1.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.43'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.126'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.51%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.87'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.46'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("E30").Value = '  +0.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0500'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.23%  '
$ws.Range("E32").Value = '  -0.27%  '
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("E34").Value = '  +0.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.42'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.902'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.136.74'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("E38").Value = '  +1.67%  '
$ws.Range("E40").Value = '  +0.72%  '
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("E42").Value = '  -1.01%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.803'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.766.18'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.15%  '
$ws.Range("E46").Value = '  -0.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.15'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0532'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.92%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.46'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.414'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.05%  '
$ws.Range("E51").Value = '  +2.76%  '
